# Edit: insert a new weekly data row for Cilantro (Terminal La Palmera de La Serena)
# immediately above the existing row 131, shifting all subsequent rows (131-229) down
# by one (to 132-230), and populate the new row 131 with the new record.
#
# The new record reuses the non-varying (boilerplate) column values from the row that
# used to be at position 131 (A, B, C, E, F, G, H, I, N, O, Q, R), but carries its own
# Date (D) and Volumen (J) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the boilerplate / unchanged values from the (soon to be shifted) row 131
# before inserting, so we can reuse them for the newly inserted row.
$colA = $ws.Range("A131").Value2
$colB = $ws.Range("B131").Value2
$colC = $ws.Range("C131").Value2
$colE = $ws.Range("E131").Value2
$colF = $ws.Range("F131").Value2
$colG = $ws.Range("G131").Value2
$colH = $ws.Range("H131").Value2
$colI = $ws.Range("I131").Value2
$colK = $ws.Range("K131").Value2
$colL = $ws.Range("L131").Value2
$colM = $ws.Range("M131").Value2
$colN = $ws.Range("N131").Value2
$colO = $ws.Range("O131").Value2
$colQ = $ws.Range("Q131").Value2
$colR = $ws.Range("R131").Value2

# Insert a new blank row at position 131; existing row 131 (and everything below)
# shifts down to row 132, etc.
$ws.Rows.Item(131).Insert()

# Populate the newly inserted row 131 with the new record.
$ws.Range("A131").Value = $colA
$ws.Range("B131").Value = $colB
$ws.Range("C131").Value = $colC
$ws.Range("D131").Value = 45062
$ws.Range("E131").Value = $colE
$ws.Range("F131").Value = $colF
$ws.Range("G131").Value = $colG
$ws.Range("H131").Value = $colH
$ws.Range("I131").Value = $colI
$ws.Range("J131").Value = 2500
$ws.Range("K131").Value = $colK
$ws.Range("L131").Value = $colL
$ws.Range("M131").Value = $colM
$ws.Range("N131").Value = $colN
$ws.Range("O131").Value = $colO
$ws.Range("P131").Value = 1500
$ws.Range("Q131").Value = $colQ
$ws.Range("R131").Value = $colR

# Make sure the new row's date cell (D) keeps the same numeric/date style as the
# other date cells in column D (style index "2" in this workbook).
$ws.Range("D131").NumberFormat = $ws.Range("D132").NumberFormat
